# Auto-generated: update cached leve-profit figures per scheduled price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3564.9092
$ws.Range("I64").Value = 2769
$ws.Range("J64").Value = 4520
$ws.Range("K64").Value = 2769
$ws.Range("L64").Value = 4520
$ws.Range("M64").Value = -2521
$ws.Range("N64").Value = -5016

$ws.Range("H67").Value = 3564.9092
$ws.Range("I67").Value = 2769
$ws.Range("J67").Value = 4520
$ws.Range("K67").Value = 2769
$ws.Range("L67").Value = 4520
$ws.Range("M67").Value = -1911
$ws.Range("N67").Value = -6236

$ws.Range("H132").Value = 1237.6444
$ws.Range("I132").Value = 1142.762
$ws.Range("K132").Value = 3428.286
$ws.Range("M132").Value = -898.2860000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4812.0786
$ws.Range("I32").Value = 3816.4092
$ws.Range("J32").Value = 11070.571
$ws.Range("K32").Value = 3816.4092
$ws.Range("L32").Value = 11070.571
$ws.Range("M32").Value = -3529.4092
$ws.Range("N32").Value = -11644.571

$ws.Range("H45").Value = 1491.7646
$ws.Range("J45").Value = 1808.091
$ws.Range("L45").Value = 1808.091
$ws.Range("N45").Value = -2562.091

$ws.Range("H61").Value = 1343.8334
$ws.Range("I61").Value = 1327.5
$ws.Range("K61").Value = 1327.5
$ws.Range("M61").Value = -1115.5

$ws.Range("H74").Value = 2525
$ws.Range("I74").Value = 997.25
$ws.Range("J74").Value = 3398
$ws.Range("K74").Value = 997.25
$ws.Range("L74").Value = 3398
$ws.Range("M74").Value = -123.25
$ws.Range("N74").Value = -5146

$ws.Range("H77").Value = 2525
$ws.Range("I77").Value = 997.25
$ws.Range("J77").Value = 3398
$ws.Range("K77").Value = 4986.25
$ws.Range("L77").Value = 16990
$ws.Range("M77").Value = -618.25
$ws.Range("N77").Value = -25726

$ws.Range("H132").Value = 1340.8611
$ws.Range("I132").Value = 1119.7037
$ws.Range("K132").Value = 3359.1111
$ws.Range("M132").Value = -829.1111000000001

$ws.Range("H136").Value = 1343.8334
$ws.Range("I136").Value = 1327.5
$ws.Range("K136").Value = 3982.5
$ws.Range("M136").Value = -1432.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 39999.5
$ws.Range("J130").Value = 39999.5
$ws.Range("L130").Value = 39999.5
$ws.Range("N130").Value = -50039.5

$ws.Range("H134").Value = 6198.2
$ws.Range("I134").Value = 6969.8096
$ws.Range("K134").Value = 20909.4288
$ws.Range("M134").Value = -18374.4288

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2673.4443
$ws.Range("I31").Value = 2472.6
$ws.Range("J31").Value = 2924.5
$ws.Range("K31").Value = 2472.6
$ws.Range("L31").Value = 2924.5
$ws.Range("M31").Value = -2177.6
$ws.Range("N31").Value = -3514.5

$ws.Range("H34").Value = 2673.4443
$ws.Range("I34").Value = 2472.6
$ws.Range("J34").Value = 2924.5
$ws.Range("K34").Value = 2472.6
$ws.Range("L34").Value = 2924.5
$ws.Range("M34").Value = -2270.6
$ws.Range("N34").Value = -3328.5

$ws.Range("H105").Value = 1064.5
$ws.Range("I105").Value = 1079.25
$ws.Range("J105").Value = 1005.5
$ws.Range("K105").Value = 1079.25
$ws.Range("L105").Value = 1005.5
$ws.Range("M105").Value = 667.75
$ws.Range("N105").Value = -4499.5

$ws.Range("H107").Value = 742.05884
$ws.Range("I107").Value = 500.2
$ws.Range("K107").Value = 500.2
$ws.Range("M107").Value = 1419.8

$ws.Range("H132").Value = 3294.95
$ws.Range("I132").Value = 2565.5715
$ws.Range("K132").Value = 7696.7145
$ws.Range("M132").Value = -5166.7145

$ws.Range("H134").Value = 1743.4062
$ws.Range("I134").Value = 1584.9259
$ws.Range("J134").Value = 2599.2
$ws.Range("K134").Value = 4754.7777
$ws.Range("L134").Value = 7797.599999999999
$ws.Range("M134").Value = -2219.7777
$ws.Range("N134").Value = -12867.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2220.8333
$ws.Range("I140").Value = 1263.5555
$ws.Range("K140").Value = 3790.6665
$ws.Range("M140").Value = 1389.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3571.2666
$ws.Range("I102").Value = 4066.5789
$ws.Range("J102").Value = 2715.7273
$ws.Range("K102").Value = 4066.5789
$ws.Range("L102").Value = 2715.7273
$ws.Range("M102").Value = -2444.5789
$ws.Range("N102").Value = -5959.7273

$ws.Range("H126").Value = 29682.865
$ws.Range("I126").Value = 2991.7693
$ws.Range("J126").Value = 44140.543
$ws.Range("K126").Value = 8975.3079
$ws.Range("L126").Value = 132421.629
$ws.Range("M126").Value = -6505.3079
$ws.Range("N126").Value = -137361.629

$ws.Range("H132").Value = 1921.4722
$ws.Range("I132").Value = 1427.75
$ws.Range("J132").Value = 2908.9167
$ws.Range("K132").Value = 4283.25
$ws.Range("L132").Value = 8726.750100000001
$ws.Range("M132").Value = -1753.25
$ws.Range("N132").Value = -13786.7501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1035.0454
$ws.Range("I22").Value = 472.14285
$ws.Range("J22").Value = 1297.7333
$ws.Range("K22").Value = 472.14285
$ws.Range("L22").Value = 1297.7333
$ws.Range("M22").Value = -177.14285
$ws.Range("N22").Value = -1887.7333

$ws.Range("H27").Value = 1035.0454
$ws.Range("I27").Value = 472.14285
$ws.Range("J27").Value = 1297.7333
$ws.Range("K27").Value = 472.14285
$ws.Range("L27").Value = 1297.7333
$ws.Range("M27").Value = -365.14285
$ws.Range("N27").Value = -1511.7333

$ws.Range("H40").Value = 23628.143
$ws.Range("I40").Value = 55001.5
$ws.Range("K40").Value = 55001.5
$ws.Range("M40").Value = -54865.5

$ws.Range("H122").Value = 5608.25
$ws.Range("I122").Value = 5883.1665
$ws.Range("J122").Value = 5333.3335
$ws.Range("K122").Value = 17649.4995
$ws.Range("L122").Value = 16000.0005
$ws.Range("M122").Value = -15199.4995
$ws.Range("N122").Value = -20900.0005

$ws.Range("H136").Value = 4081.6428
$ws.Range("I136").Value = 3063.7368
$ws.Range("J136").Value = 6230.5557
$ws.Range("K136").Value = 9191.2104
$ws.Range("L136").Value = 18691.6671
$ws.Range("M136").Value = -6641.2104
$ws.Range("N136").Value = -23791.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 113053.86
$ws.Range("I122").Value = 196069.25
$ws.Range("K122").Value = 588207.75
$ws.Range("M122").Value = -585757.75

$ws.Range("H126").Value = 19640.572
$ws.Range("I126").Value = 27876
$ws.Range("J126").Value = 8660
$ws.Range("K126").Value = 83628
$ws.Range("L126").Value = 25980
$ws.Range("M126").Value = -81158
$ws.Range("N126").Value = -30920

